$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.775.75"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "2.082.70"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'234.52"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'58.82"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "'0.0788"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").Value = "2.390.63"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "'14.83"
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").Value = "'21.13"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "'0.774"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "'5.32"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "2.081.07"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "37.706.75"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").Value = "'228.93"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "'169.31"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  +4.50%  "
$ws.Range("D28").Value = "'9.01"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").Value = "'1.84"
$ws.Range("E36").Value = "  +2.91%  "
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'5.42"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("D40").Value = "'0.0981"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("D41").Value = "'98.58"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").Value = "'2.87"
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("D44").Value = "1.462.63"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").Value = "'4.31"
$ws.Range("E45").Value = "  +3.56%  "
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'16.45"
$ws.Range("E47").Value = "  +5.84%  "
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("D50").Value = "'3.03"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "2.276.07"
$ws.Range("E51").Value = "  -0.30%  "
